$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 472 (shifts old rows 472-560 down to 476-564)
$ws.Rows.Item(472).Resize(4).Insert()

# Common (unchanged across these 4 new rows) field values
$mercado = 'Mercado Mayorista Lo Valledor de Santiago'
$region = 'Metropolitana'
$tipo = 'Fruta'
$producto = 'Otros'
$categoria = 'Chirimoya'
$variedad = 'Cultivar IV Región'
$unidad = '$/bandeja 10 kilos'
$origen = 'Provincia de Limarí'

# Data for the 4 newly inserted rows (472-475)
$newRows = @(
    @{ Row = 472; Fecha = 45244; Calidad = 'Especial'; Volumen = 200; PMin = 21000; PMax = 21000; PProm = 21000; PKg = 2100; KgUnidad = 10 },
    @{ Row = 473; Fecha = 45244; Calidad = 'Primera';  Volumen = 300; PMin = 17000; PMax = 17000; PProm = 17000; PKg = 1700; KgUnidad = 10 },
    @{ Row = 474; Fecha = 45244; Calidad = 'Segunda';  Volumen = 275; PMin = 13000; PMax = 13000; PProm = 13000; PKg = 1300; KgUnidad = 10 },
    @{ Row = 475; Fecha = 45244; Calidad = 'Tercera';  Volumen = 250; PMin = 10000; PMax = 10000; PProm = 10000; PKg = 1000; KgUnidad = 10 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 6
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = 100107002
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
